$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.06606273387611693, 0.7819567165648904, 1.95711398432533, 1.398968900413919, 1.411041781561629)
    3  = @(1.090855880406542, 1.636952682927139, 9.587480958581839, 3.096365766278564, 2.926681515795593)
    4  = @(1.146478367687022, 2.300800508439017, 16.92217468558993, 4.113657093826602, 3.990775337436762)
    5  = @(1.174361413305189, 2.149695694135142, 16.96319945630839, 4.118640486411553, 3.988576556960211)
    6  = @(0.7872502277514816, 2.189782397188602, 17.30793470436396, 4.160280604041506, 4.128345827086147)
    7  = @(0.8319901022676859, 1.971575915675119, 16.97606752739902, 4.120202364860131, 4.088078109409235)
    8  = @(0.8848395731913589, 2.077403147304855, 17.89474323381673, 4.230217870726841, 4.192169283686689)
    9  = @(0.7851577424002278, 2.363787703475875, 28.62182591829283, 5.349937001338692, 5.422694984943732)
    10 = @(-0.5578088075490664, 1.742949889236707, 10.92979066599325, 3.306023391628264, 3.381635471449998)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
